$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 126, shifting existing rows 126-250 down to 127-251.
$ws.Rows(126).Insert()

# Populate the new row 126 with a fresh weekly data point.
$ws.Range("A126").Value = 8
$ws.Range("B126").Value = "Terminal La Palmera de La Serena"
$ws.Range("C126").Value = "Coquimbo"
$ws.Range("D126").Value = 45167
$ws.Range("E126").Value = 4
$ws.Range("F126").Value = 100112044
$ws.Range("G126").Value = "Perejil"
$ws.Range("H126").Value = "Sin especificar"
$ws.Range("I126").Value = "Primera"
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 2000
$ws.Range("L126").Value = 2500
$ws.Range("M126").Value = 2250
$ws.Range("N126").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O126").Value = "Provincia del Elquí"
$ws.Range("P126").Value = 1500
$ws.Range("Q126").Value = 1.5
$ws.Range("R126").Value = "Hortaliza"
